$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: the sentence "...Python, Natural language processing
# library, SQL, API's." was split across three runs ("...p" / "r" /
# "ocessing..."). Re-running Find/Replace over the full span coalesces
# it back into a single run with the same (already correct) text.
# ---------------------------------------------------------------------
$find1 = $d.Content.Find
$find1.Text = "Python, Natural language processing"
$find1.Replacement.Text = "Python, Natural language processing"
$find1.Execute($find1.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)

# ---------------------------------------------------------------------
# Change 2: the two blank paragraphs right after "Exceptional Investing
# Performance." (between it and "This document with Github links:")
# gain explicit Bold = False / BoldCs = False on their (empty) run.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $t = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Exceptional Investing Performance.") {
        $blank1 = $d.Paragraphs($i + 1)
        $blank1.Range.Font.Bold = $false
        $blank1.Range.Font.BoldBi = $false
        $blank2 = $d.Paragraphs($i + 2)
        $blank2.Range.Font.Bold = $false
        $blank2.Range.Font.BoldBi = $false
        break
    }
}

# ---------------------------------------------------------------------
# Change 3: drop the stray "." run that trails the Github document
# hyperlink at the very end of the document.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $t = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "https://github.com/GGTechGuru/PROFESSIONAL_DOCUMENTS/blob/main/GERARD_GOLD_SWQA_INVESTING_RESEARCH_TUTORING.docx.") {
        $e = $para.Range.End
        $periodRange = $d.Range($e - 2, $e - 1)
        $periodRange.Delete()
        break
    }
}

Write-Output "edits applied"
